# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in the title cell
# - Update case data for several countries (Belgica, Banglades, Chequia)
# - Indonesia overtakes Sudafrica in total cases, so the two countries swap
#   places in the ranking (row 34 becomes Indonesia, row 35 becomes Sudafrica)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Title / timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 11:05"

# Belgica (row 21) - updated case numbers
$ws.Range("B21").Value = 57342
$ws.Range("C21").Value = 250
$ws.Range("D21").Value = 15297
$ws.Range("E21").Value = 32733
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 32
$ws.Range("H21").Value = 9312

# Banglades (row 27) - updated case numbers
$ws.Range("B27").Value = 35585
$ws.Range("C27").Value = 1975
$ws.Range("D27").Value = 7334
$ws.Range("E27").Value = 27750
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 21
$ws.Range("H27").Value = 501

# Indonesia overtakes Sudafrica (22750 > 22583): swap the two rows.
# Row 34 now holds Indonesia with new data, row 35 now holds Sudafrica
# with the data that previously sat in row 34.
$ws.Range("A34").Value = "Indonesia"
$ws.Range("B34").Value = 22750
$ws.Range("C34").Value = 479
$ws.Range("D34").Value = 5642
$ws.Range("E34").Value = 15717
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 19
$ws.Range("H34").Value = 1391

$ws.Range("A35").Value = "Sudafrica"
$ws.Range("B35").Value = 22583
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 11100
$ws.Range("E35").Value = 11054
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 429

# Chequia (row 54) - updated case numbers
$ws.Range("B54").Value = 8957
$ws.Range("C54").Value = 2
$ws.Range("D54").Value = 6083
$ws.Range("E54").Value = 2559
